# 🔄 MAJ automatique BRVM via GitHub Actions
# Update the "Top_YTD" sheet: re-rank the rows by the new (huge) values and
# refresh the "Progression YTD (%)" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top_YTD")

# New ordering (A) + refreshed values (B) for rows 2..11
$titles = @(
    "BRVM - SERVICES PUBLICS",
    "BRVM - AUTRES SECTEURS",
    "VIVO ENERGY CI",
    "CFAO MOTORS CI",
    "NEI-CEDA CI",
    "SUCRIVOIRE",
    "SETAO CI",
    "SAFCA CI",
    "AIR LIQUIDE CI",
    "UNIWAX CI"
)

$values = @(
    [double]"3.545696671258368e+123",
    [double]"1.615967677106988e+88",
    [double]"1.029917107421899e+87",
    [double]"1.883572412868977e+84",
    [double]"6.092003722527686e+82",
    [double]"7.00539501208747e+80",
    [double]"5.520884726404072e+79",
    [double]"1.758205686436238e+79",
    [double]"6.124509488580116e+77",
    [double]"2.702158929032325e+72"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $titles[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
